# Fix: SO Ref Num (and related order data) was repeating the same values
# for rows 17-21; re-populate rows 17-23 with the corrected, de-duplicated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17
$ws.Range('A17').Value = ' '
$ws.Range('B17').Value = ' '
$ws.Range('C17').Value = ' '
$ws.Range('D17').Value = ' '
$ws.Range('E17').Value = ' '
$ws.Range('F17').Value = ' '
$ws.Range('G17').Value = ' '
$ws.Range('H17').Value = ' '
$ws.Range('I17').Value = ' '
$ws.Range('J17').Value = ' '
$ws.Range('K17').Value = ' '
$ws.Range('L17').Value = ' '
$ws.Range('M17').Value = ' '
$ws.Range('N17').Value = ' '
$ws.Range('O17').Value = ' '
$ws.Range('P17').Value = ' '
$ws.Range('Q17').Value = ' '
$ws.Range('R17').Value = ' '
$ws.Range('S17').Value = ' '
$ws.Range('T17').Value = ' '
$ws.Range('U17').Value = ' '
$ws.Range('V17').Value = ' '
$ws.Range('W17').Value = ' '
$ws.Range('X17').Value = ' '
$ws.Range('Y17').Value = ' '
$ws.Range('Z17').Value = ' '
$ws.Range('AA17').Value = ' '
$ws.Range('AB17').Value = ' '
$ws.Range('AC17').Value = ' '
$ws.Range('AD17').Value = ' '

# Row 18
$c = $ws.Range('A18')
$c.NumberFormat = '@'
$c.Value = '03/11/2024'
$c.Style = 'Normal'
$ws.Range('B18').Value = 'SO240311001'
$ws.Range('C18').Value = 'NO'
$ws.Range('D18').Value = 'ab'
$ws.Range('E18').Value = 'abc@abc.com'
$ws.Range('F18').Value = '(789)456-1230 '
$ws.Range('G18').Value = 'YES'
$ws.Range('H18').Value = 'YES'
$ws.Range('I18').Value = 'artist'
$ws.Range('J18').Value = 'title'
$ws.Range('K18').Value = 5
$ws.Range('L18').Value = 15
$ws.Range('M18').Value = ' '
$ws.Range('N18').Value = 'AMS'
$ws.Range('O18').Value = 'LP'
$ws.Range('P18').Value = 'abake'
$ws.Range('Q18').Value = 'PICKUP'
$ws.Range('R18').Value = ''
$ws.Range('S18').Value = ''
$ws.Range('T18').Value = ''
$ws.Range('U18').Value = ''
$ws.Range('V18').Value = 'NO'
$ws.Range('W18').Value = ''
$ws.Range('X18').Value = ''
$ws.Range('Y18').Value = 'NO'
$ws.Range('Z18').Value = ''
$ws.Range('AA18').Value = ''
$ws.Range('AB18').Value = 'NO'
$ws.Range('AC18').Value = ''
$ws.Range('AD18').Value = ''

# Row 19
$c = $ws.Range('A19')
$c.NumberFormat = '@'
$c.Value = '03/11/2024'
$c.Style = 'Normal'
$ws.Range('B19').Value = 'SO240311002'
$ws.Range('C19').Value = 'NO'
$ws.Range('D19').Value = 'ab'
$ws.Range('E19').Value = 'asdsf'
$ws.Range('F19').Value = '(789)465-1230'
$ws.Range('G19').Value = 'YES'
$ws.Range('H19').Value = 'YES'
$ws.Range('I19').Value = 'artist'
$ws.Range('J19').Value = 'title'
$ws.Range('K19').Value = 10
$ws.Range('L19').Value = 50.99
$ws.Range('M19').Value = ' '
$ws.Range('N19').Value = 'AMS'
$ws.Range('O19').Value = 'LP'
$ws.Range('P19').Value = 'abake'
$ws.Range('Q19').Value = 'SHIPPING'
$ws.Range('R19').Value = '2120 Septiembre dr'
$ws.Range('S19').Value = 'ELPaso'
$ws.Range('T19').Value = 'TX'
$c = $ws.Range('U19')
$c.NumberFormat = '@'
$c.Value = '79935'
$c.Style = 'Normal'
$ws.Range('V19').Value = 'NO'
$ws.Range('W19').Value = ''
$ws.Range('X19').Value = ''
$ws.Range('Y19').Value = 'NO'
$ws.Range('Z19').Value = ''
$ws.Range('AA19').Value = ''
$ws.Range('AB19').Value = 'NO'
$ws.Range('AC19').Value = ''
$ws.Range('AD19').Value = ''

# Row 20
$c = $ws.Range('A20')
$c.NumberFormat = '@'
$c.Value = '03/11/2024'
$c.Style = 'Normal'
$ws.Range('B20').Value = 'SO240311003'
$ws.Range('C20').Value = 'NO'
$ws.Range('D20').Value = 'ab'
$ws.Range('E20').Value = 'ab'
$ws.Range('F20').Value = '(789)123-4567'
$ws.Range('G20').Value = 'YES'
$ws.Range('H20').Value = 'YES'
$ws.Range('I20').Value = 'artist'
$ws.Range('J20').Value = 'title'
$ws.Range('K20').Value = 5
$ws.Range('L20').Value = 90
$ws.Range('M20').Value = ' '
$ws.Range('N20').Value = 'AMS'
$ws.Range('O20').Value = 'OTHER'
$ws.Range('P20').Value = 'abake'
$ws.Range('Q20').Value = 'PICKUP'
$ws.Range('R20').Value = ''
$ws.Range('S20').Value = ''
$ws.Range('T20').Value = ''
$ws.Range('U20').Value = ''
$ws.Range('V20').Value = 'NO'
$ws.Range('W20').Value = ''
$ws.Range('X20').Value = ''
$ws.Range('Y20').Value = 'NO'
$ws.Range('Z20').Value = ''
$ws.Range('AA20').Value = ''
$ws.Range('AB20').Value = 'NO'
$ws.Range('AC20').Value = ''
$ws.Range('AD20').Value = ''

# Row 21
$c = $ws.Range('A21')
$c.NumberFormat = '@'
$c.Value = '03/11/2024'
$c.Style = 'Normal'
$ws.Range('B21').Value = 'SO240311004'
$ws.Range('C21').Value = 'NO'
$ws.Range('D21').Value = 'Ashley'
$ws.Range('E21').Value = 'ab'
$ws.Range('F21').Value = '(159)753-2856'
$ws.Range('G21').Value = 'YES'
$ws.Range('H21').Value = 'YES'
$ws.Range('I21').Value = 'artist'
$ws.Range('J21').Value = 'title'
$ws.Range('K21').Value = 60
$ws.Range('L21').Value = 822
$ws.Range('M21').Value = ' '
$ws.Range('N21').Value = 'AEC'
$ws.Range('O21').Value = 'DVD'
$ws.Range('P21').Value = 'abake'
$ws.Range('Q21').Value = 'PICKUP'
$ws.Range('R21').Value = ''
$ws.Range('S21').Value = ''
$ws.Range('T21').Value = ''
$ws.Range('U21').Value = ''
$ws.Range('V21').Value = 'NO'
$ws.Range('W21').Value = ''
$ws.Range('X21').Value = ''
$ws.Range('Y21').Value = 'NO'
$ws.Range('Z21').Value = ''
$ws.Range('AA21').Value = ''
$ws.Range('AB21').Value = 'NO'
$ws.Range('AC21').Value = ''
$ws.Range('AD21').Value = ''

# Row 22
$c = $ws.Range('A22')
$c.NumberFormat = '@'
$c.Value = '03/11/2024'
$c.Style = 'Normal'
$ws.Range('B22').Value = 'SO240311005'
$ws.Range('C22').Value = 'NO'
$ws.Range('D22').Value = 'ab'
$ws.Range('E22').Value = 'ababa'
$ws.Range('F22').Value = '(789)456-1230'
$ws.Range('G22').Value = 'YES'
$ws.Range('H22').Value = 'YES'
$ws.Range('I22').Value = 'artist'
$ws.Range('J22').Value = 'title`'
$ws.Range('K22').Value = 9
$ws.Range('L22').Value = 55
$ws.Range('M22').Value = ' '
$ws.Range('N22').Value = 'AEC'
$ws.Range('O22').Value = 'LP'
$ws.Range('P22').Value = 'abake'
$ws.Range('Q22').Value = 'PICKUP'
$ws.Range('R22').Value = ''
$ws.Range('S22').Value = ''
$ws.Range('T22').Value = ''
$ws.Range('U22').Value = ''
$ws.Range('V22').Value = 'NO'
$ws.Range('W22').Value = ''
$ws.Range('X22').Value = ''
$ws.Range('Y22').Value = 'NO'
$ws.Range('Z22').Value = ''
$ws.Range('AA22').Value = ''
$ws.Range('AB22').Value = 'NO'
$ws.Range('AC22').Value = ''
$ws.Range('AD22').Value = ''

# Row 23
$c = $ws.Range('A23')
$c.NumberFormat = '@'
$c.Value = '03/11/2024'
$c.Style = 'Normal'
$ws.Range('B23').Value = 'SO240311006'
$ws.Range('C23').Value = 'NO'
$ws.Range('D23').Value = 'ab'
$ws.Range('E23').Value = 'ababababa'
$ws.Range('F23').Value = '(528)525-8225'
$ws.Range('G23').Value = 'YES'
$ws.Range('H23').Value = 'YES'
$ws.Range('I23').Value = 'art'
$ws.Range('J23').Value = 'title'
$ws.Range('K23').Value = 16
$ws.Range('L23').Value = 98
$ws.Range('M23').Value = ' '
$ws.Range('N23').Value = 'me'
$ws.Range('O23').Value = 'BLU-RAY'
$ws.Range('P23').Value = 'abake'
$ws.Range('Q23').Value = 'PICKUP'
$ws.Range('R23').Value = 'N/A'
$ws.Range('S23').Value = 'N/A'
$ws.Range('T23').Value = 'N/A'
$ws.Range('U23').Value = 'N/A'
$ws.Range('V23').Value = 'NO'
$ws.Range('W23').Value = ''
$ws.Range('X23').Value = ''
$ws.Range('Y23').Value = 'NO'
$ws.Range('Z23').Value = ''
$ws.Range('AA23').Value = ''
$ws.Range('AB23').Value = 'NO'
$ws.Range('AC23').Value = ''
$ws.Range('AD23').Value = ''

Write-Host "Row 17-23 data corrected; SO Ref Num duplication resolved."